# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TB/d2S/K/IP values (columns B-E), Win (F) unchanged, sum (G) = B+C+D+E
$data = @(
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027),
    @(0.2881169905109251, 0.04103571897497393, 3.223369029078222, 0.5333859586016987, 0, 4.085907697165819),
    @(0.2881169905109251, 0.3048912486333797, 0.7210945179870265, 13.86384647080068, 1, 15.17794922793202),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $row++
}

$wb.Save()
